# Cost profile sheet refactor: rename header labels, collapse the
# baseline columns down to a single "bl_one" column (drop baseline
# two / baseline three), and refresh "last" (column C) values to
# match the (now single) baseline column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update column headers in row 1 ---
$ws.Range("B1").Value = "current"
$ws.Range("C1").Value = "last"
$ws.Range("D1").Value = "bl_one"

# --- 2. Refresh "last" (column C) values for rows 9-22 so that they
#        match the recomputed baseline (column D) numbers. ---
$ws.Range("C9").Value  = 848.6500000000001
$ws.Range("C10").Value = 1194.22
$ws.Range("C11").Value = 2728.757082855212
$ws.Range("C12").Value = 3505.88
$ws.Range("C13").Value = 612.01
$ws.Range("C14").Value = 217.61
$ws.Range("C15").Value = 219.97
$ws.Range("C16").Value = 223.06
$ws.Range("C17").Value = 226.19
$ws.Range("C18").Value = 229.36
$ws.Range("C19").Value = 232.56
$ws.Range("C20").Value = 216.62
$ws.Range("C21").Value = 145.26
$ws.Range("C22").Value = 51.19

# Baseline recompute also nudged D9 by a floating point ULP.
$ws.Range("D9").Value = 848.6500000000001

# --- 3. Drop the now-redundant "Baseline two" / "Baseline three"
#        columns (E and F) entirely. ---
$ws.Columns("E:F").Delete()
